{"js": "// Locate the paragraph that ends the \"May, 2011 by A. Nagy\" byline (the\n// right-aligned date/author line right under the title block) so the new\n// \"Last updated\" line can be inserted immediately after it.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nlet bylinePara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Nagy\") >= 0) {\n    bylinePara = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!bylinePara) {\n  throw new Error(\"Could not locate the 'Nagy' byline paragraph\");\n}\n\n// The stray leftover \"_GoBack\" bookmark further down the document (next to\n// the myDialog.Show() snippet) is removed \u2014 Word re-creates \"_GoBack\" at the\n// location of the most recent edit, which is the new paragraph below.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Insert a new right-aligned paragraph after the byline containing a\n// \"Last updated, Date : \" label followed by a DATE field (complex field\n// form: begin/instrText/separate/cached-result/end) and a fresh \"_GoBack\"\n// bookmark wrapping the field, matching what Word stamps after an edit.\nconst afterByline = bylinePara.getRange(\"After\");\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  \"<w:p>\" +\n  \"<w:pPr>\" +\n  '<w:jc w:val=\"right\"/>' +\n  \"</w:pPr>\" +\n  \"<w:r>\" +\n  '<w:t xml:space=\"preserve\">Last updated, Date : </w:t>' +\n  \"</w:r>\" +\n  \"<w:r>\" +\n  '<w:fldChar w:fldCharType=\"begin\"/>' +\n  \"</w:r>\" +\n  \"<w:r>\" +\n  '<w:instrText xml:space=\"preserve\"> DATE \\\\@ \"MMMM d, yyyy\" </w:instrText>' +\n  \"</w:r>\" +\n  \"<w:r>\" +\n  '<w:fldChar w:fldCharType=\"separate\"/>' +\n  \"</w:r>\" +\n  \"<w:r>\" +\n  \"<w:rPr>\" +\n  \"<w:noProof/>\" +\n  \"</w:rPr>\" +\n  \"<w:t>March 19, 2014</w:t>\" +\n  \"</w:r>\" +\n  \"<w:r>\" +\n  '<w:fldChar w:fldCharType=\"end\"/>' +\n  \"</w:r>\" +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  \"</w:p>\" +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nafterByline.insertOoxml(ooxml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The stray leftover \"_GoBack\" bookmark further down the document (next to\n# the myDialog.Show() code snippet) is removed first \u2014 Word re-creates\n# \"_GoBack\" at the location of the most recent edit, which will be the new\n# paragraph inserted below.\n$d.Bookmarks(\"_GoBack\").Delete()\n\n# Locate the paragraph that ends the \"May, 2011 by A. Nagy\" byline (the\n# right-aligned date/author line right under the title block) so the new\n# \"Last updated\" line can be inserted immediately after it.\n$target = $null\n$targetIdx = 0\n$idx = 0\nforeach ($p in $d.Paragraphs) {\n    $idx = $idx + 1\n    if ($p.Range.Text -like \"*Nagy*\") {\n        $target = $p\n        $targetIdx = $idx\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not locate the 'Nagy' byline paragraph\"\n}\n\n# Create a new, empty paragraph right after the byline.\n$target.Range.InsertParagraphAfter()\n$newPara = $d.Paragraphs.Item($targetIdx + 1)\n\n# Fill it in via raw OOXML: a right-aligned \"Last updated, Date : \" label\n# followed by a DATE field (complex field form: begin/instrText/separate/\n# cached-result/end) and a fresh \"_GoBack\" bookmark wrapping the field,\n# matching what Word stamps after an edit.\n$xml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:pPr><w:jc w:val=\"right\"/></w:pPr><w:r><w:t xml:space=\"preserve\">Last updated, Date : </w:t></w:r><w:r><w:fldChar w:fldCharType=\"begin\"/></w:r><w:r><w:instrText xml:space=\"preserve\"> DATE \\@ \"MMMM d, yyyy\" </w:instrText></w:r><w:r><w:fldChar w:fldCharType=\"separate\"/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t>March 19, 2014</w:t></w:r><w:r><w:fldChar w:fldCharType=\"end\"/></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>'\n$newPara.Range.InsertXML($xml)\n"}
